# Generate Report for Handback
# Updates the handoff/handback timestamps and the zh-cn priority value
# to reflect a newer report-generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 2 & 4
$wsOverview.Range("G2").Value = "2016-08-30 10:16:22"
$wsOverview.Range("G4").Value = "2016-08-30 10:16:22"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K) for rows 2 & 4
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-30 10:16:18"
$wsZhCn.Range("H4").Value = "2016-08-30 10:16:18"
$wsZhCn.Range("K2").Value = "2016-08-30 10:16:35"
$wsZhCn.Range("K4").Value = "2016-08-30 10:16:35"

# de-de sheet: Priority (E), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K) for rows 2 & 4
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-30 10:16:22"
$wsDeDe.Range("H4").Value = "2016-08-30 10:16:22"
$wsDeDe.Range("K2").Value = "2016-08-30 10:16:42"
$wsDeDe.Range("K4").Value = "2016-08-30 10:16:42"
